$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the "Description" column (column D) and
# shift the existing "Description" data over to column E.
$ws.Range("D1").EntireColumn.Insert()

# New header for the inserted column.
$ws.Range("D1").Value = "Price"

# New price values for each expense row.
$ws.Range("D2").Value = 101.25
$ws.Range("D3").Value = 114.23
$ws.Range("D4").Value = 318.75
